$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "2025/12/03 18:00"
$ws.Range("B26").Value = "-"
$ws.Range("C26").Value = "-"
$ws.Range("D26").Value = "-"
$ws.Range("E26").Value = "-"
$ws.Range("F26").Value = "-"
$ws.Range("G26").Value = "-"
